$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (label "0.0")
$ws.Range("B2").Value = 0.949438202247191
$ws.Range("C2").Value = 0.8644501278772379
$ws.Range("D2").Value = 0.9049531459170014
$ws.Range("E2").Value = 782

# Row 3 (label "1.0")
$ws.Range("B3").Value = 0.6293706293706294
$ws.Range("D3").Value = 0.7171314741035857

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.8577154308617234
$ws.Range("C4").Value = 0.8577154308617234
$ws.Range("D4").Value = 0.8577154308617234
$ws.Range("E4").Value = 0.8577154308617234

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.7894044158089102
$ws.Range("C5").Value = 0.8488917306052857
$ws.Range("D5").Value = 0.8110423100102935
$ws.Range("E5").Value = 998

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.8801650602218029
$ws.Range("C6").Value = 0.8577154308617234
$ws.Range("D6").Value = 0.8643023632399495
$ws.Range("E6").Value = 998
